$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new data points on row 7
$ws.Range("C7").Value = 17109
$ws.Range("E7").Value = 16411.8

# Update the active selection to E11
$ws.Range("E11").Select()
